# B6-PowerPoint.pptx edit
#
# 1) Re-style the three financial-statement tables (slides 14-16) from the
#    "No Style, Table Grid" built-in table style to "Medium Style 2 -
#    Accent 1" (the default themed table style).
# 2) Re-colour the theme actually used by the deck (the slide master's
#    theme, physically stored as ppt/theme/theme2.xml) from the custom
#    "Integral / Red Violet" palette back to the stock "Office" palette,
#    using the real ThemeColorScheme/ThemeColor object model (the same
#    mechanism PowerPoint uses under Design > Colors).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Tables: swap the table style id on every table shape on slides
#    14, 15 and 16.
# ---------------------------------------------------------------------
$newTableStyleId = "{3C38210F-F05C-41DA-8E48-7C7468327579}"

for ($slideIdx = 14; $slideIdx -le 16; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours: restore the stock "Office" colour scheme.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$colors = $slide1.ThemeColorScheme

# index : Office colour (was the Integral / Red Violet colour)
$colors.Colors(1).RGB  = 0         # dk1      000000 (unchanged)
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF (unchanged)
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hyperlink        0563C1
$colors.Colors(12).RGB = 7491477   # followed hyperlink 954F72
